$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 9 (shifts existing rows 9..101 down to 10..102,
# carrying their formatting/styles with them - matches dimension change A1:R101 -> A1:R102).
$ws.Rows.Item(9).EntireRow.Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44761
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112001
$ws.Range("G9").Value = "Berenjena"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 220
$ws.Range("K9").Value = 11000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 11455
$ws.Range("N9").Value = "$/caja 60 unidades"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 191
$ws.Range("Q9").Value = 60
$ws.Range("R9").Value = "Hortaliza"
